# Applies the cryptos-list data refresh described by the commit:
# "Updated cryptos list on Tue Aug  1 07:25:51 UTC 2023 with GitHub Actions"
#
# Every data cell on the sheet is stored as an inline/shared STRING (t="inlineStr"
# in the source OOXML), even when the text looks like a plain number (e.g. "244.28").
# Excel.Range.Value auto-converts number-looking text to a real numeric value, which
# would corrupt those cells (e.g. "0.9996" -> 0.9996, "64.11" -> 64.11 as a Double,
# dropping the literal formatting). To keep them as text we pre-format any cell whose
# new value parses as a number with a Text ("@") number format before writing it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cell -> new value, in the same order as the diff
$updates = [ordered]@{
    'D2' = '28.902.22'
    'E2' = '  -1.75%  '
    'D3' = '1.829.90'
    'E3' = '  -2.05%  '
    'E4' = '  -0.09%  '
    'D5' = '244.85'
    'E5' = '  +0.35%  '
    'D6' = '0.6890'
    'E6' = '  -2.60%  '
    'E7' = '  -0.06%  '
    'D8' = '0.07650'
    'E8' = '  -2.82%  '
    'D9' = '0.3053'
    'E9' = '  -2.63%  '
    'D10' = '23.51'
    'E10' = '  -4.29%  '
    'D11' = '0.07813'
    'E11' = '  -1.89%  '
    'D12' = '1.834.73'
    'E12' = '  -2.82%  '
    'D13' = '5.072'
    'E13' = '  -2.67%  '
    'D14' = '90.48'
    'E14' = '  -3.18%  '
    'D15' = '0.6773'
    'E15' = '  -3.42%  '
    'D16' = '6.431'
    'E16' = '  -1.51%  '
    'D17' = '0.000008292'
    'E17' = '  -0.83%  '
    'D18' = '28.918.86'
    'E18' = '  -1.85%  '
    'D19' = '242.77'
    'E19' = '  -3.92%  '
    'D20' = '2.081.01'
    'E20' = '  -2.56%  '
    'D21' = '12.67'
    'E21' = '  -3.48%  '
    'D22' = '0.9996'
    'D23' = '7.437'
    'E23' = '  -2.66%  '
    'E24' = '  -0.09%  '
    'D25' = '0.1471'
    'E25' = '  -5.39%  '
    'D26' = '161.59'
    'D27' = '8.790'
    'E27' = '  -2.36%  '
    'D28' = '18.17'
    'E28' = '  -2.92%  '
    'D29' = '1.553'
    'E29' = '  +3.46%  '
    'D30' = '4.210'
    'E30' = '  -2.82%  '
    'D31' = '4.129'
    'E31' = '  -3.00%  '
    'E32' = '  -2.57%  '
    'D33' = '0.05107'
    'E33' = '  -3.99%  '
    'D34' = '0.7566'
    'E34' = '  +0.93%  '
    'D35' = '1.831'
    'E35' = '  -3.27%  '
    'E36' = '  -2.52%  '
    'E37' = '  -1.22%  '
    'D38' = '0.01844'
    'E38' = '  -2.48%  '
    'D39' = '1.232.98'
    'E39' = '  -3.22%  '
    'D40' = '2.689'
    'E40' = '  -2.21%  '
    'D41' = '0.9276'
    'E41' = '  +3.70%  '
    'D42' = '108.69'
    'E42' = '  -0.45%  '
    'D43' = '0.9994'
    'E43' = '  -0.10%  '
    'D44' = '5.674'
    'E44' = '  -6.63%  '
    'D45' = '9.544'
    'E45' = '  -0.19%  '
    'D46' = '1.980.27'
    'E46' = '  -2.29%  '
    'E47' = '  -0.30%  '
    'B48' = 'BabyDogeCoin'
    'C48' = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
    'D48' = '0.00000000121'
    'E48' = '  -5.01%  '
    'B49' = 'Aave'
    'C49' = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
    'D49' = '64.11'
    'E49' = '  -10.13%  '
    'D50' = '1.735'
    'E50' = '  -3.45%  '
    'D51' = '0.4193'
    'E51' = '  -2.69%  '
}

foreach ($cellRef in $updates.Keys) {
    $value = $updates[$cellRef]
    $range = $ws.Range($cellRef)
    if ($value -match '^[+-]?\d+(\.\d+)?$') {
        # Force text storage so the numeric-looking string is not coerced to a Double
        $range.NumberFormat = "@"
    }
    $range.Value = $value
}
